# "Fruta / hortaliza, semanal" - weekly data refresh.
# A new daily price record for Acelga (Macroferia Regional de Talca) is
# inserted at the top of this sheet's data block (row 433), pushing the
# existing rows 433-510 down by one (to 434-511).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 433; this shifts rows
# 433..510 down to 434..511 (and the sheet dimension grows to R511).
$ws.Rows(433).Insert()

# Populate the newly inserted row with the new record.
$ws.Cells.Item(433, 1).Value = 5
$ws.Cells.Item(433, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(433, 3).Value = "Maule"
$ws.Cells.Item(433, 4).Value = 45209
$ws.Cells.Item(433, 5).Value = 7
$ws.Cells.Item(433, 6).Value = 100112009
$ws.Cells.Item(433, 7).Value = "Acelga"
$ws.Cells.Item(433, 8).Value = "Sin especificar"
$ws.Cells.Item(433, 9).Value = "Primera"
$ws.Cells.Item(433, 10).Value = 500
$ws.Cells.Item(433, 11).Value = 1800
$ws.Cells.Item(433, 12).Value = 1800
$ws.Cells.Item(433, 13).Value = 1800
$ws.Cells.Item(433, 14).Value = "$/docena de atados (4 kilos)"
$ws.Cells.Item(433, 15).Value = "Provincia de Curicó"
$ws.Cells.Item(433, 16).Value = 450
$ws.Cells.Item(433, 17).Value = 4
$ws.Cells.Item(433, 18).Value = "Hortaliza"
